$d = $word.ActiveDocument

# The location line lives in the second paragraph (contact-info block).
# Scope the Find/Replace to that paragraph's range so the other two
# "College Station, TX" occurrences (Education / Work sections) are left
# untouched.
$p = $d.Paragraphs.Item(2)
$rng = $p.Range
$rng.Find.Execute("College Station, TX", $true, $true, $false, $false, $false, `
                   $true, 1, $false, "Dallas, TX", 2)
